$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WF")

$ws.Range("D8").Value = 7524600
$ws.Range("E8").Value = 7490800
$ws.Range("F8").Value = 7654400
$ws.Range("G8").Value = 8105900
$ws.Range("H8").Value = 8354200
$ws.Range("I8").Value = 9584300
$ws.Range("J8").Value = 9763900

$ws.Range("D15").Value = -165000
$ws.Range("E15").Value = -221800
$ws.Range("F15").Value = -211900
$ws.Range("G15").Value = -200600
$ws.Range("H15").Value = -208000
$ws.Range("I15").Value = -196400
$ws.Range("J15").Value = -177900

$ws.Range("D17").Value = 3689200
$ws.Range("E17").Value = 3857500
$ws.Range("F17").Value = 4442700
$ws.Range("G17").Value = 5134300
$ws.Range("H17").Value = 6337500
$ws.Range("I17").Value = 6909700
$ws.Range("J17").Value = 7124100

$ws.Range("D18").Value = 3835400
$ws.Range("E18").Value = 3633300
$ws.Range("F18").Value = 3211800
$ws.Range("G18").Value = 2971600
$ws.Range("H18").Value = 2016700
$ws.Range("I18").Value = 2674600
$ws.Range("J18").Value = 2639700

$ws.Range("D20").Value = -2119800
$ws.Range("E20").Value = -2266300
$ws.Range("F20").Value = -1934100
$ws.Range("G20").Value = -2237400
$ws.Range("H20").Value = -1763500
$ws.Range("I20").Value = -1233200
$ws.Range("J20").Value = -578900

$ws.Range("D21").Value = 1923100
$ws.Range("E21").Value = 1588800
$ws.Range("F21").Value = 1489600
$ws.Range("G21").Value = 951800
$ws.Range("H21").Value = 517500
$ws.Range("I21").Value = 1697500
$ws.Range("J21").Value = 2301600

$ws.Range("D23").Value = 1715600
$ws.Range("E23").Value = 1367000
$ws.Range("F23").Value = 1277700
$ws.Range("G23").Value = 734300
$ws.Range("H23").Value = 253100
$ws.Range("I23").Value = 1441400
$ws.Range("J23").Value = 2060800

$ws.Range("D24").Value = 369100
$ws.Range("E24").Value = 242800
$ws.Range("F24").Value = 331400
$ws.Range("G24").Value = 253600
$ws.Range("H24").Value = 30900
$ws.Range("I24").Value = 314000
$ws.Range("J24").Value = 491900

$ws.Range("D26").Value = 1346500
$ws.Range("E26").Value = 1124200
$ws.Range("F26").Value = 946300
$ws.Range("G26").Value = 480700
$ws.Range("H26").Value = 222300
$ws.Range("I26").Value = 1127400
$ws.Range("J26").Value = 1568900

$ws.Range("D27").Value = 1183700
$ws.Range("E27").Value = 928200
$ws.Range("F27").Value = 770700
$ws.Range("G27").Value = 338900
$ws.Range("H27").Value = 116700
$ws.Range("I27").Value = 1000600
$ws.Range("J27").Value = 1437400

$ws.Range("G29").Value = 685200
$ws.Range("H29").Value = -615700
$ws.Range("I29").Value = 412600
$ws.Range("J29").Value = 456200

$ws.Range("D32").Value = 2119800
$ws.Range("E32").Value = 2266300
$ws.Range("F32").Value = 1934100
$ws.Range("G32").Value = 2237400
$ws.Range("H32").Value = 1763500
$ws.Range("I32").Value = 1233200
$ws.Range("J32").Value = 578900

$ws.Range("D33").Value = 1183700
$ws.Range("E33").Value = 928200
$ws.Range("F33").Value = 770700
$ws.Range("G33").Value = 1024200
$ws.Range("H33").Value = -499000
$ws.Range("I33").Value = 1413300
$ws.Range("J33").Value = 1893700

$ws.Range("D35").Value = 1183700
$ws.Range("E35").Value = 928200
$ws.Range("F35").Value = 770700
$ws.Range("G35").Value = 1024200
$ws.Range("H35").Value = -499000
$ws.Range("I35").Value = 1413300
$ws.Range("J35").Value = 1893700

$ws.Range("D41").Value = 5323900
$ws.Range("E41").Value = 6243600
$ws.Range("F41").Value = 5309900
$ws.Range("G41").Value = 4577300
$ws.Range("H41").Value = 9056100
$ws.Range("I41").Value = 5085000
$ws.Range("J41").Value = 5647100

$ws.Range("D42").Value = 20932300
$ws.Range("E42").Value = 13888400
$ws.Range("F42").Value = 12461100
$ws.Range("G42").Value = 15487500
$ws.Range("H42").Value = 38324100
$ws.Range("I42").Value = 49687900

$ws.Range("D47").Value = 367000
$ws.Range("E47").Value = 386300
$ws.Range("F47").Value = 566600
$ws.Range("G47").Value = 570600
$ws.Range("H47").Value = 1567600
$ws.Range("I47").Value = 913400
$ws.Range("J47").Value = 816800

$ws.Range("D48").Value = 2507000
$ws.Range("E48").Value = 2478500
$ws.Range("F48").Value = 2484000
$ws.Range("G48").Value = 2515600
$ws.Range("H48").Value = 5063600
$ws.Range("I48").Value = 5205200
$ws.Range("J48").Value = 3197500

$ws.Range("D49").Value = 456400
$ws.Range("E49").Value = 425700
$ws.Range("F49").Value = 369400
$ws.Range("G49").Value = 260200
$ws.Range("H49").Value = 473300
$ws.Range("I49").Value = 762800
$ws.Range("J49").Value = 394100

$ws.Range("D52").Value = 289300
$ws.Range("E52").Value = 268700
$ws.Range("F52").Value = 201100
$ws.Range("G52").Value = 234000
$ws.Range("H52").Value = 75071700
$ws.Range("I52").Value = 210100
$ws.Range("J52").Value = 70400

$ws.Range("D54").Value = 278340000
$ws.Range("E54").Value = 273401000
$ws.Range("F54").Value = 256836000
$ws.Range("G54").Value = 237738000
$ws.Range("H54").Value = 299808000
$ws.Range("I54").Value = 287850000
$ws.Range("J54").Value = 275257000

$ws.Range("D57").Value = 5801500
$ws.Range("E57").Value = 13033100
$ws.Range("F57").Value = 7371600
$ws.Range("G57").Value = 7298400
$ws.Range("H57").Value = 13514900
$ws.Range("I57").Value = 1048000

$ws.Range("H58").Value = 6864200

$ws.Range("D59").Value = 2008600
$ws.Range("E59").Value = 1959900
$ws.Range("F59").Value = 1768900
$ws.Range("G59").Value = 2325000
$ws.Range("H59").Value = 4051000
$ws.Range("I59").Value = 3245300
$ws.Range("J59").Value = 241300

$ws.Range("D61").Value = 36941800
$ws.Range("E61").Value = 33771200
$ws.Range("F61").Value = 34482400
$ws.Range("G61").Value = 34714100
$ws.Range("H61").Value = 30283200
$ws.Range("I61").Value = 46588000

$ws.Range("D62").Value = 419200
$ws.Range("E62").Value = 453300
$ws.Range("F62").Value = 559400
$ws.Range("G62").Value = 694600
$ws.Range("H62").Value = 1417700
$ws.Range("I62").Value = 1931100
$ws.Range("J62").Value = 1119700

$ws.Range("D66").Value = 260418000
$ws.Range("E66").Value = 255461000
$ws.Range("F66").Value = 239950000
$ws.Range("G66").Value = 221913000
$ws.Range("H66").Value = 284133000
$ws.Range("I66").Value = 271397000
$ws.Range("J66").Value = 259836000

$ws.Range("D72").Value = 13745600
$ws.Range("E72").Value = 12858200
$ws.Range("F72").Value = 12079000
$ws.Range("G72").Value = 12465500
$ws.Range("H72").Value = 23078300
$ws.Range("I72").Value = 12215600
$ws.Range("J72").Value = 10932000

$ws.Range("D76").Value = 17922000
$ws.Range("E76").Value = 17939800
$ws.Range("F76").Value = 16885900
$ws.Range("G76").Value = 15825500
$ws.Range("H76").Value = 15674700
$ws.Range("I76").Value = 16452400
$ws.Range("J76").Value = 15421100

$ws.Range("D81").Value = 1183700
$ws.Range("E81").Value = 928200
$ws.Range("F81").Value = 770700
$ws.Range("G81").Value = 1024200
$ws.Range("H81").Value = -499000
$ws.Range("I81").Value = 1413300
$ws.Range("J81").Value = 1893700

$ws.Range("D83").Value = 207500
$ws.Range("E83").Value = 221800
$ws.Range("F83").Value = 211900
$ws.Range("G83").Value = 217600
$ws.Range("H83").Value = 264400
$ws.Range("I83").Value = 256100
$ws.Range("J83").Value = 240800

$ws.Range("D89").Value = -1741500
$ws.Range("E89").Value = 4316700
$ws.Range("F89").Value = -337800
$ws.Range("G89").Value = 247500
$ws.Range("H89").Value = -3513600
$ws.Range("I89").Value = -126800
$ws.Range("J89").Value = -942600

$ws.Range("D91").Value = -142800
$ws.Range("E91").Value = -115300
$ws.Range("F91").Value = -113900
$ws.Range("G91").Value = -123800
$ws.Range("H91").Value = -140300
$ws.Range("I91").Value = -222200
$ws.Range("J91").Value = -204600

$ws.Range("D94").Value = 1657700
$ws.Range("E94").Value = -3718500
$ws.Range("F94").Value = 1121200
$ws.Range("G94").Value = -1004100
$ws.Range("H94").Value = -345800
$ws.Range("I94").Value = 2069900
$ws.Range("J94").Value = 2124000

$ws.Range("D96").Value = -296200
$ws.Range("E96").Value = -148100
$ws.Range("F96").Value = -444400
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = -333000
$ws.Range("I96").Value = -302100
$ws.Range("J96").Value = -177300

$ws.Range("D100").Value = -53800
$ws.Range("E100").Value = 293100
$ws.Range("F100").Value = -423600
$ws.Range("G100").Value = 274100
$ws.Range("H100").Value = 3678400
$ws.Range("I100").Value = -2355300
$ws.Range("J100").Value = 154300

$ws.Range("D101").Value = -463400
$ws.Range("E101").Value = -57600
$ws.Range("F101").Value = 239700
$ws.Range("G101").Value = 34100
$ws.Range("H101").Value = -83600
$ws.Range("I101").Value = -150700
$ws.Range("J101").Value = 24700

$ws.Range("D102").Value = -601100
$ws.Range("E102").Value = 833600
$ws.Range("F102").Value = 599500
$ws.Range("G102").Value = -448400
$ws.Range("H102").Value = -264700
$ws.Range("I102").Value = -562800
$ws.Range("J102").Value = 1360500
